$d = $word.ActiveDocument

# Paragraph 7 (0-based) / Word's Paragraphs collection is 1-based -> index 8
$p1 = $d.Paragraphs.Item(8).Range
$p1.Find.Execute("{{c.name}} | {{c.phone}} | {{c.relation}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{INS `$c.name}} | {{INS `$c.phone}} | {{INS `$c.relation}}", 2)

# Paragraph 8 (0-based) -> Word index 9
$p2 = $d.Paragraphs.Item(9).Range
$p2.Find.Execute("{{END-FOR}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{END-FOR c}}", 2)

# Paragraph 11 (0-based) -> Word index 12
$p3 = $d.Paragraphs.Item(12).Range
$p3.Find.Execute("{{m.name}} | {{m.dosage}} | {{m.schedule}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{INS `$m.name}} | {{INS `$m.dosage}} | {{INS `$m.schedule}}", 2)

# Paragraph 12 (0-based) -> Word index 13
$p4 = $d.Paragraphs.Item(13).Range
$p4.Find.Execute("{{END-FOR}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{END-FOR m}}", 2)

# Paragraph 18 (0-based) -> Word index 19
$p5 = $d.Paragraphs.Item(19).Range
$p5.Find.Execute("{{p}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{INS `$p}}", 2)

# Paragraph 19 (0-based) -> Word index 20
$p6 = $d.Paragraphs.Item(20).Range
$p6.Find.Execute("{{END-FOR}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{END-FOR p}}", 2)
